# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The sheet holds metadata about the dataset columns, laid out as:
#   Row 1: column names (ano, comarca-codigo, municipio-nombre, diputados, ...)
#   Row 2: sdmx/iaest role for each column (dimension vs measure)
#   Row 3: "dim" / "medida" classification
#   Row 4: data type / concept scheme for each column
#   Row 5: mapping workbook file used for curation
#
# "municipio-nombre" (column C) and "diputados" (column D) get their
# curated dimension/measure roles swapped, and the now-unused
# mapping-diputados.xlsx curation file reference is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# municipio-nombre (column C) becomes a curated dimension (refArea), with
# a proper URI concept reference instead of being treated as a measure.
$ws.Cells.Item(2, 3).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(3, 3).Value = "dim"
$ws.Cells.Item(4, 3).Value = "URI-Municipio"

# diputados (column D) becomes a curated measure instead of a dimension.
$ws.Cells.Item(2, 4).Value = "iaest-measure:diputados"
$ws.Cells.Item(3, 4).Value = "medida"
$ws.Cells.Item(4, 4).Value = "xsd:int"

# Its mapping workbook reference is no longer needed now that it isn't a
# curated dimension requiring a mapping file.
$ws.Cells.Item(5, 4).Clear()
